$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = -0.2254024683979639

$ws.Range("J11").Value = -0.1253231084953424
$ws.Range("K11").Value = -0.3352267436446591

$ws.Range("I12").Value = 0.2284633975843539
$ws.Range("J12").Value = 0.01855976243503714

$ws.Range("H13").Value = 0.08028600715190851
$ws.Range("I13").Value = -0.1296176279974082

$ws.Range("G14").Value = -0.07715998185224648
$ws.Range("H14").Value = -0.2870636170015632

$ws.Range("F15").Value = 0.4234994746738243
$ws.Range("G15").Value = 0.2135958395245076

$ws.Range("E16").Value = 0.1431415941383551
$ws.Range("F16").Value = -0.06676204101096155

$ws.Range("D17").Value = 0.3151164519833668
$ws.Range("E17").Value = 0.1052128168340501

$ws.Range("C18").Value = 0.009253912237035311
$ws.Range("D18").Value = -0.2006497229122814

$ws.Range("B19").Value = 0.6215838649243215
$ws.Range("C19").Value = 0.4116802297750048

$ws.Range("B20").Value = -0.2766911554241067
